# Update the header date and regenerate all 26 division problems/answers
# in the worksheet table (each cell's text is unique, so a straightforward
# Find/Replace-All per pair is unambiguous and order-independent).
$d = $word.ActiveDocument

# Header date
$d.Content.Find.Execute("2025-05-20 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-21 Wednesday", 2) | Out-Null

# Table answers (three-digit ÷ one-digit problems)
$d.Content.Find.Execute("288÷3=96, 0", $true, $false, $false, $false, $false, $true, 1, $false, "585÷7=83, 4", 2) | Out-Null
$d.Content.Find.Execute("724÷8=90, 4", $true, $false, $false, $false, $false, $true, 1, $false, "419÷8=52, 3", 2) | Out-Null
$d.Content.Find.Execute("611÷3=203, 2", $true, $false, $false, $false, $false, $true, 1, $false, "997÷7=142, 3", 2) | Out-Null
$d.Content.Find.Execute("956÷6=159, 2", $true, $false, $false, $false, $false, $true, 1, $false, "790÷4=197, 2", 2) | Out-Null
$d.Content.Find.Execute("527÷6=87, 5", $true, $false, $false, $false, $false, $true, 1, $false, "268÷3=89, 1", 2) | Out-Null
$d.Content.Find.Execute("256÷2=128, 0", $true, $false, $false, $false, $false, $true, 1, $false, "126÷5=25, 1", 2) | Out-Null
$d.Content.Find.Execute("804÷9=89, 3", $true, $false, $false, $false, $false, $true, 1, $false, "635÷6=105, 5", 2) | Out-Null
$d.Content.Find.Execute("518÷4=129, 2", $true, $false, $false, $false, $false, $true, 1, $false, "782÷7=111, 5", 2) | Out-Null
$d.Content.Find.Execute("782÷5=156, 2", $true, $false, $false, $false, $false, $true, 1, $false, "431÷9=47, 8", 2) | Out-Null
$d.Content.Find.Execute("377÷7=53, 6", $true, $false, $false, $false, $false, $true, 1, $false, "562÷3=187, 1", 2) | Out-Null
$d.Content.Find.Execute("625÷8=78, 1", $true, $false, $false, $false, $false, $true, 1, $false, "882÷3=294, 0", 2) | Out-Null
$d.Content.Find.Execute("888÷7=126, 6", $true, $false, $false, $false, $false, $true, 1, $false, "984÷8=123, 0", 2) | Out-Null
$d.Content.Find.Execute("298÷4=74, 2", $true, $false, $false, $false, $false, $true, 1, $false, "504÷5=100, 4", 2) | Out-Null
$d.Content.Find.Execute("909÷2=454, 1", $true, $false, $false, $false, $false, $true, 1, $false, "577÷5=115, 2", 2) | Out-Null
$d.Content.Find.Execute("640÷4=160, 0", $true, $false, $false, $false, $false, $true, 1, $false, "675÷7=96, 3", 2) | Out-Null
$d.Content.Find.Execute("611÷5=122, 1", $true, $false, $false, $false, $false, $true, 1, $false, "827÷8=103, 3", 2) | Out-Null
$d.Content.Find.Execute("279÷3=93, 0", $true, $false, $false, $false, $false, $true, 1, $false, "868÷3=289, 1", 2) | Out-Null
$d.Content.Find.Execute("759÷2=379, 1", $true, $false, $false, $false, $false, $true, 1, $false, "745÷2=372, 1", 2) | Out-Null
$d.Content.Find.Execute("914÷9=101, 5", $true, $false, $false, $false, $false, $true, 1, $false, "752÷5=150, 2", 2) | Out-Null
$d.Content.Find.Execute("393÷4=98, 1", $true, $false, $false, $false, $false, $true, 1, $false, "642÷8=80, 2", 2) | Out-Null
$d.Content.Find.Execute("109÷7=15, 4", $true, $false, $false, $false, $false, $true, 1, $false, "744÷6=124, 0", 2) | Out-Null
$d.Content.Find.Execute("219÷9=24, 3", $true, $false, $false, $false, $false, $true, 1, $false, "706÷2=353, 0", 2) | Out-Null
$d.Content.Find.Execute("546÷3=182, 0", $true, $false, $false, $false, $false, $true, 1, $false, "109÷2=54, 1", 2) | Out-Null
$d.Content.Find.Execute("678÷2=339, 0", $true, $false, $false, $false, $false, $true, 1, $false, "119÷7=17, 0", 2) | Out-Null
$d.Content.Find.Execute("250÷6=41, 4", $true, $false, $false, $false, $false, $true, 1, $false, "477÷9=53, 0", 2) | Out-Null

Write-Output "Applied 26 replacements."
